$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "deadVolume" row of experiment-condition data (row 6, between
# the existing Feed conc row and the results header block).
$ws.Range("A6").Value = "deadVolume"
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = "mL"

# Match B6's fill to the rest of the yellow "value" cells above it (B2:B5)
# before the alignment pass so they all end up sharing one style.
$ws.Range("B6").Interior.Color = 65535

# Right-align the whole block of numeric input cells (B2:B6).
$ws.Range("B2:B6").HorizontalAlignment = -4152

# Column F is no longer used anywhere on the sheet - clear out the two
# leftover formatted-but-empty cells so the sheet shrinks back to A:E.
$ws.Range("F7:F8").Clear()

# Leave the sheet scrolled back to the top with the new cell selected,
# like a user would after typing the new row in.
$ws.Range("B6").Select()
